$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.580.53"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.922.51"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "'325.63"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "'0.4814"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").Value = "'0.4058"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").Value = "'0.08217"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").Value = "'23.66"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "1.939.49"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").Value = "'6.071"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").Value = "'7.262"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").Value = "'91.62"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "'0.06857"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "'1.012"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "'0.00001040"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "'17.62"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "29.583.09"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "'5.686"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").Value = "'11.92"
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").Value = "'2.183"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "2.139.15"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").Value = "'155.75"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").Value = "'6.455"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "'2.093"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").Value = "'120.64"
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("D31").Value = "'1.015"
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("D32").Value = "'0.09624"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").Value = "'5.616"
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("D34").Value = "'3.556"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").Value = "'1.378"
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("D36").Value = "'0.06349"
$ws.Range("E36").Value = "  +4.08%  "
$ws.Range("D37").Value = "'0.02290"
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("D38").Value = "'1.183"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").Value = "'0.5956"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").Value = "'10.76"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").Value = "'7.871"
$ws.Range("E42").Value = "  -1.36%  "
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").Value = "'2.453"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").Value = "'1.287"
$ws.Range("E45").Value = "  +3.08%  "
$ws.Range("D46").Value = "'12.40"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").Value = "'0.07499"
$ws.Range("E47").Value = "  -2.69%  "
$ws.Range("D48").Value = "'0.5567"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").Value = "'119.27"
$ws.Range("E50").Value = "  +3.14%  "
$ws.Range("D51").Value = "'2.440"
$ws.Range("E51").Value = "  +3.63%  "
